$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Cells.Item(1,1).Value = "Datos actualizados a 27 de Mayo de 2020 a las 16:05"

# Update country rows (re-sorted ranking + refreshed case counts)
# Row 4
$ws.Cells.Item(4,2).Value = 1729771
$ws.Cells.Item(4,3).Value = 4496
$ws.Cells.Item(4,4).Value = 480273
$ws.Cells.Item(4,5).Value = 1148753
$ws.Cells.Item(4,7).Value = 173
$ws.Cells.Item(4,8).Value = 100745
# Row 11
$ws.Cells.Item(11,2).Value = 181530
$ws.Cells.Item(11,3).Value = 242
$ws.Cells.Item(11,5).Value = 10232
# Row 13
$ws.Cells.Item(13,2).Value = 154181
$ws.Cells.Item(13,3).Value = 3388
$ws.Cells.Item(13,4).Value = 65394
$ws.Cells.Item(13,5).Value = 84414
$ws.Cells.Item(13,7).Value = 29
$ws.Cells.Item(13,8).Value = 4373
# Row 47
$ws.Cells.Item(47,4).Value = 4349
$ws.Cells.Item(47,5).Value = 8387
$ws.Cells.Item(47,7).Value = 8
$ws.Cells.Item(47,8).Value = 492
# Row 57
$ws.Cells.Item(57,2).Value = 8391
$ws.Cells.Item(57,3).Value = 8
$ws.Cells.Item(57,5).Value = 429
# Row 67
$ws.Cells.Item(67,5).Value = 1279
$ws.Cells.Item(67,7).Value = 1
$ws.Cells.Item(67,8).Value = 313
# Row 111
$ws.Cells.Item(111,1).Value = "Zambia"
$ws.Cells.Item(111,2).Value = 1057
$ws.Cells.Item(111,3).Value = 137
$ws.Cells.Item(111,4).Value = 779
$ws.Cells.Item(111,5).Value = 271
$ws.Cells.Item(111,8).Value = 7
# Row 112
$ws.Cells.Item(112,1).Value = "Tunez"
$ws.Cells.Item(112,2).Value = 1051
$ws.Cells.Item(112,3).Value = 0
$ws.Cells.Item(112,4).Value = 929
$ws.Cells.Item(112,5).Value = 74
$ws.Cells.Item(112,8).Value = 48
# Row 113
$ws.Cells.Item(113,1).Value = "Albania"
$ws.Cells.Item(113,2).Value = 1050
$ws.Cells.Item(113,3).Value = 21
$ws.Cells.Item(113,4).Value = 812
$ws.Cells.Item(113,5).Value = 205
$ws.Cells.Item(113,8).Value = 33
# Row 114
$ws.Cells.Item(114,1).Value = "Guinea Ecuatorial"
$ws.Cells.Item(114,2).Value = 1043
$ws.Cells.Item(114,4).Value = 165
$ws.Cells.Item(114,5).Value = 866
$ws.Cells.Item(114,8).Value = 12
# Row 115
$ws.Cells.Item(115,1).Value = "Costa Rica"
$ws.Cells.Item(115,2).Value = 956
$ws.Cells.Item(115,4).Value = 634
$ws.Cells.Item(115,5).Value = 312
$ws.Cells.Item(115,8).Value = 10
# Row 116
$ws.Cells.Item(116,1).Value = "Niger"
$ws.Cells.Item(116,2).Value = 952
$ws.Cells.Item(116,4).Value = 796
$ws.Cells.Item(116,5).Value = 93
$ws.Cells.Item(116,8).Value = 63
# Row 117
$ws.Cells.Item(117,1).Value = "Republica de Chipre"
$ws.Cells.Item(117,2).Value = 939
$ws.Cells.Item(117,4).Value = 594
$ws.Cells.Item(117,5).Value = 328
$ws.Cells.Item(117,8).Value = 17
# Row 157
$ws.Cells.Item(157,4).Value = 126
$ws.Cells.Item(157,5).Value = 74
# Row 199
$ws.Cells.Item(199,1).Value = "Belice"
$ws.Cells.Item(199,4).Value = 16
$ws.Cells.Item(199,8).Value = 2
# Row 200
$ws.Cells.Item(200,1).Value = "Nueva Caledonia"
$ws.Cells.Item(200,4).Value = 18
$ws.Cells.Item(200,8).Value = 0
# Row 201
$ws.Cells.Item(201,1).Value = "Santa Lucia"
# Row 207
$ws.Cells.Item(207,1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(207,4).Value = 10
$ws.Cells.Item(207,8).Value = 1
# Row 208
$ws.Cells.Item(208,1).Value = "Groenlandia"
$ws.Cells.Item(208,4).Value = 11
$ws.Cells.Item(208,8).Value = 0
# Row 210
$ws.Cells.Item(210,1).Value = "Seychelles"
$ws.Cells.Item(210,4).Value = 11
$ws.Cells.Item(210,8).Value = 0
# Row 211
$ws.Cells.Item(211,1).Value = "Montserrat"
$ws.Cells.Item(211,4).Value = 10
$ws.Cells.Item(211,8).Value = 1
